# Apply the "Add data for 2021-11-20" update:
#  - Rename the sheet to reflect the new "through" date (11-11 -> 11-12)
#  - Update the month label for November to match
#  - Update October's 2021 (H) value
#  - Update November's values for all year columns
#  - Update the Total row values for all year columns

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename sheet (tab name) and the label cell for November
$ws.Name = "Through 2021-11-12"
$ws.Range("A12").Value = "November (through 11-12)"

# October (row 11) - only 2021 column (H) changes
$ws.Range("H11").Value = 194

# November (row 12) - all year columns change
$ws.Range("B12").Value = 14
$ws.Range("C12").Value = 28
$ws.Range("D12").Value = 49
$ws.Range("E12").Value = 26
$ws.Range("F12").Value = 21
$ws.Range("G12").Value = 77
$ws.Range("H12").Value = 83

# Total (row 13) - all year columns change
$ws.Range("B13").Value = 272
$ws.Range("C13").Value = 514
$ws.Range("D13").Value = 759
$ws.Range("E13").Value = 641
$ws.Range("F13").Value = 503
$ws.Range("G13").Value = 1134
$ws.Range("H13").Value = 1525
